$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number: force text storage so Excel
# does not auto-convert the string into a numeric value (matches original
# inlineStr/shared-string cell type in the source workbook).
$forceTextCells = @(
    @{ Addr = "D5"; Val = '595.92' }
    @{ Addr = "D6"; Val = '150.17' }
    @{ Addr = "D10"; Val = '5.70' }
    @{ Addr = "D11"; Val = '0.381' }
    @{ Addr = "D13"; Val = '27.71' }
    @{ Addr = "D19"; Val = '4.67' }
    @{ Addr = "D20"; Val = '347.87' }
    @{ Addr = "D21"; Val = '6.87' }
    @{ Addr = "D22"; Val = '0.999' }
    @{ Addr = "D23"; Val = '5.72' }
    @{ Addr = "D24"; Val = '66.39' }
    @{ Addr = "D28"; Val = '578.44' }
    @{ Addr = "D29"; Val = '8.20' }
    @{ Addr = "D35"; Val = '5.26' }
    @{ Addr = "D36"; Val = '169.15' }
    @{ Addr = "D38"; Val = '0.999' }
    @{ Addr = "D40"; Val = '19.36' }
    @{ Addr = "D42"; Val = '169.11' }
    @{ Addr = "D43"; Val = '39.91' }
    @{ Addr = "D46"; Val = '21.39' }
    @{ Addr = "D47"; Val = '0.628' }
    @{ Addr = "D48"; Val = '0.0249' }
    @{ Addr = "D49"; Val = '1.98' }
    @{ Addr = "D50"; Val = '0.0965' }
    @{ Addr = "D51"; Val = '19.33' }
)
foreach ($item in $forceTextCells) {
    $ws.Range($item.Addr).NumberFormat = "@"
}
foreach ($item in $forceTextCells) {
    $ws.Range($item.Addr).Value = $item.Val
}
foreach ($item in $forceTextCells) {
    $ws.Range($item.Addr).Style = "Normal"
}

# Remaining cells already round-trip as text (percentages always contain a
# "%" sign/spaces; a few D-column prices use "." as a thousands separator,
# which is not a valid numeric literal), so a plain assignment is enough.
$ws.Range("D2").Value = '63.670.06'
$ws.Range("D3").Value = '2.626.73'
$ws.Range("E3").Value = '  -0.82%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  -1.37%  '
$ws.Range("E6").Value = '  +2.27%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("E9").Value = '  +0.52%  '
$ws.Range("E10").Value = '  +1.91%  '
$ws.Range("E11").Value = '  +3.24%  '
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").Value = '3.098.23'
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").Value = '63.511.56'
$ws.Range("E15").Value = '  +0.05%  '
$ws.Range("E16").Value = '  +2.19%  '
$ws.Range("D17").Value = '2.635.22'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("E18").Value = '  +7.20%  '
$ws.Range("E19").Value = '  +2.28%  '
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("E21").Value = '  -0.82%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("E23").Value = '  +2.31%  '
$ws.Range("E24").Value = '  -0.71%  '
$ws.Range("E25").Value = '  +11.55%  '
$ws.Range("E26").Value = '  -0.92%  '
$ws.Range("E27").Value = '  +1.44%  '
$ws.Range("E28").Value = '  +0.79%  '
$ws.Range("E29").Value = '  +2.52%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("D33").Value = '0.0₃0844'
$ws.Range("E33").Value = '  +2.73%  '
$ws.Range("E34").Value = '  -0.42%  '
$ws.Range("E35").Value = '  +0.76%  '
$ws.Range("E36").Value = '  +0.30%  '
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  +0.23%  '
$ws.Range("E40").Value = '  +1.34%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("E44").Value = '  +3.95%  '
$ws.Range("E45").Value = '  +5.32%  '
$ws.Range("E46").Value = '  -3.83%  '
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("E48").Value = '  +1.10%  '
$ws.Range("E49").Value = '  +4.78%  '
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("E51").Value = '  +2.54%  '
